$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $c = $ws.Range($CellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $NewValue
    $c.Style = $origStyle
}

Set-TextValue "D2" "300.77"
Set-TextValue "E2" "0.66%"
Set-TextValue "D3" "32.16"
Set-TextValue "E3" "1.73%"
Set-TextValue "D4" "4.960"
Set-TextValue "E4" "-3.82%"
Set-TextValue "D5" "0.07881"
Set-TextValue "E5" "-2.59%"
Set-TextValue "D6" "2.073"
Set-TextValue "E6" "-17.17%"
Set-TextValue "D7" "7.798"
Set-TextValue "E7" "-0.12%"
Set-TextValue "D8" "3.833"
Set-TextValue "E8" "-2.20%"
Set-TextValue "D9" "0.9257"
Set-TextValue "E9" "-0.15%"
Set-TextValue "D10" "0.1739"
Set-TextValue "E10" "-1.11%"
Set-TextValue "D11" "0.07953"
Set-TextValue "E11" "5.93%"
Set-TextValue "D12" "0.08656"
Set-TextValue "E12" "-2.08%"
Set-TextValue "D13" "0.03098"
Set-TextValue "E13" "3.22%"
Set-TextValue "D14" "0.1001"
Set-TextValue "E14" "0.06%"
Set-TextValue "D15" "0.001517"
Set-TextValue "E15" "0.69%"
Set-TextValue "D16" "0.005893"
Set-TextValue "E16" "-1.22%"
Set-TextValue "E17" "2,098.44%"
Set-TextValue "D18" "3.458"
Set-TextValue "E18" "-2.03%"
Set-TextValue "E19" "-2.10%"
Set-TextValue "D20" "0.3276"
Set-TextValue "E20" "0.12%"
Set-TextValue "D22" "4.292"
Set-TextValue "E22" "2.89%"
Set-TextValue "E23" "6.73%"
Set-TextValue "D24" "0.04593"
Set-TextValue "E24" "-0.74%"
Set-TextValue "D25" "0.001230"
Set-TextValue "E25" "-0.97%"
Set-TextValue "D26" "0.004421"
Set-TextValue "E26" "-2.25%"
Set-TextValue "D27" "0.0001251"
Set-TextValue "E27" "4.22%"
Set-TextValue "D39" "0.01712"
Set-TextValue "E39" "-2.06%"
Set-TextValue "D40" "0.04763"
Set-TextValue "E40" "3.66%"
Set-TextValue "D41" "0.007435"
Set-TextValue "E41" "7.71%"
Set-TextValue "D42" "0.1353"
Set-TextValue "E42" "-1.44%"
Set-TextValue "D43" "0.002352"
Set-TextValue "E43" "7.36%"
Set-TextValue "D44" "0.01121"
Set-TextValue "E44" "8.70%"
Set-TextValue "D45" "0.00005991"
Set-TextValue "E45" "-5.36%"
Set-TextValue "E46" "0.05%"
Set-TextValue "D47" "0.003392"
Set-TextValue "E47" "-59.64%"
Set-TextValue "D48" "0.8205"
Set-TextValue "E48" "9.62%"
Set-TextValue "D49" "0.00002102"
Set-TextValue "E49" "0.05%"
Set-TextValue "D50" "0.0002002"
Set-TextValue "E50" "0.05%"
